# Refresh the cryptocurrency price / link / 1h-volume table with the latest
# scraped data (scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns B/C/E hold plain text (names, links, padded percentages)
# and are safe to assign directly. Column D holds price strings that can
# look numeric (e.g. "1.00", "7.10", "46.512.32"); to keep them stored as
# text -- exactly as they were before the edit -- and to avoid float
# round-off / trailing-zero loss, they are written with a temporary "@"
# (text) number format, then the style is reset to "Normal" so no stray
# formatting is left behind on the cell.

$priceUpdates = @(
    @{ Cell = "D2"; Value = "46.512.32" }
    @{ Cell = "D3"; Value = "2.248.26" }
    @{ Cell = "D5"; Value = "298.84" }
    @{ Cell = "D6"; Value = "98.96" }
    @{ Cell = "D7"; Value = "0.558" }
    @{ Cell = "D8"; Value = "1.00" }
    @{ Cell = "D9"; Value = "0.507" }
    @{ Cell = "D10"; Value = "35.35" }
    @{ Cell = "D11"; Value = "0.0775" }
    @{ Cell = "D12"; Value = "7.10" }
    @{ Cell = "D14"; Value = "2.591.61" }
    @{ Cell = "D15"; Value = "2.251.50" }
    @{ Cell = "D16"; Value = "13.46" }
    @{ Cell = "D17"; Value = "46.534.75" }
    @{ Cell = "D18"; Value = "0.787" }
    @{ Cell = "D19"; Value = "12.75" }
    @{ Cell = "D20"; Value = "0.0₃0917" }
    @{ Cell = "D21"; Value = "5.84" }
    @{ Cell = "D22"; Value = "65.02" }
    @{ Cell = "D23"; Value = "247.50" }
    @{ Cell = "D24"; Value = "2.80" }
    @{ Cell = "D25"; Value = "0.999" }
    @{ Cell = "D26"; Value = "1.85" }
    @{ Cell = "D27"; Value = "42.13" }
    @{ Cell = "D28"; Value = "2.23" }
    @{ Cell = "D29"; Value = "9.62" }
    @{ Cell = "D30"; Value = "19.71" }
    @{ Cell = "D32"; Value = "145.39" }
    @{ Cell = "D33"; Value = "5.37" }
    @{ Cell = "D34"; Value = "0.0763" }
    @{ Cell = "D35"; Value = "3.13" }
    @{ Cell = "D36"; Value = "0.113" }
    @{ Cell = "D38"; Value = "15.86" }
    @{ Cell = "D39"; Value = "1.69" }
    @{ Cell = "D40"; Value = "3.81" }
    @{ Cell = "D41"; Value = "0.0295" }
    @{ Cell = "D42"; Value = "3.17" }
    @{ Cell = "D43"; Value = "0.999" }
    @{ Cell = "D44"; Value = "1.95" }
    @{ Cell = "D45"; Value = "1.812.32" }
    @{ Cell = "D46"; Value = "89.78" }
    @{ Cell = "D47"; Value = "0.186" }
    @{ Cell = "D48"; Value = "70.96" }
    @{ Cell = "D49"; Value = "4.78" }
    @{ Cell = "D50"; Value = "93.20" }
    @{ Cell = "D51"; Value = "2.469.03" }
)

foreach ($u in $priceUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

$textUpdates = @(
    @{ Cell = "E2"; Value = "  +3.91%  " }
    @{ Cell = "E3"; Value = "  -0.54%  " }
    @{ Cell = "E4"; Value = "  +0.02%  " }
    @{ Cell = "E5"; Value = "  -0.89%  " }
    @{ Cell = "E6"; Value = "  +5.06%  " }
    @{ Cell = "E7"; Value = "  -0.93%  " }
    @{ Cell = "E8"; Value = "  +0.14%  " }
    @{ Cell = "E9"; Value = "  -0.12%  " }
    @{ Cell = "E10"; Value = "  +3.92%  " }
    @{ Cell = "E11"; Value = "  -1.60%  " }
    @{ Cell = "E12"; Value = "  -1.12%  " }
    @{ Cell = "E13"; Value = "  -1.14%  " }
    @{ Cell = "E14"; Value = "  -0.43%  " }
    @{ Cell = "E15"; Value = "  -0.19%  " }
    @{ Cell = "E16"; Value = "  -0.68%  " }
    @{ Cell = "E18"; Value = "  -1.10%  " }
    @{ Cell = "E19"; Value = "  -1.00%  " }
    @{ Cell = "E20"; Value = "  -0.44%  " }
    @{ Cell = "E21"; Value = "  -3.47%  " }
    @{ Cell = "E22"; Value = "  -0.82%  " }
    @{ Cell = "E23"; Value = "  +4.05%  " }
    @{ Cell = "E24"; Value = "  -2.68%  " }
    @{ Cell = "E25"; Value = "  +0.14%  " }
    @{ Cell = "E26"; Value = "  -1.52%  " }
    @{ Cell = "E27"; Value = "  +2.06%  " }
    @{ Cell = "E28"; Value = "  -2.79%  " }
    @{ Cell = "E29"; Value = "  +0.44%  " }
    @{ Cell = "E30"; Value = "  +1.00%  " }
    @{ Cell = "E31"; Value = "  +9.03%  " }
    @{ Cell = "E32"; Value = "  -4.47%  " }
    @{ Cell = "E33"; Value = "  -2.95%  " }
    @{ Cell = "E34"; Value = "  -3.31%  " }
    @{ Cell = "E35"; Value = "  +6.18%  " }
    @{ Cell = "E36"; Value = "  +8.59%  " }
    @{ Cell = "E37"; Value = "  -1.87%  " }
    @{ Cell = "E38"; Value = "  +16.89%  " }
    @{ Cell = "E39"; Value = "  -2.89%  " }
    @{ Cell = "E41"; Value = "  -4.53%  " }
    @{ Cell = "E42"; Value = "  -1.83%  " }
    @{ Cell = "E43"; Value = "  -0.14%  " }
    @{ Cell = "E44"; Value = "  +2.96%  " }
    @{ Cell = "E45"; Value = "  +4.34%  " }
    @{ Cell = "E46"; Value = "  +18.48%  " }
    @{ Cell = "E47"; Value = "  -3.90%  " }
    @{ Cell = "E48"; Value = "  +2.35%  " }
    @{ Cell = "E49"; Value = "  +2.04%  " }
    @{ Cell = "E50"; Value = "  -2.82%  " }
    @{ Cell = "B51"; Value = "RocketPoolETH" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth" }
    @{ Cell = "E51"; Value = "  -0.62%  " }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

